$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Share" column (K) for several scene rows so that a new
# group gets created when a payer logs in with a single-clone scene.
$ws.Range("K10").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("K15").Value = 0

# Move the active cell selection to K10 (was K15).
$ws.Range("K10").Select()
